# Updated CHE_grids model - 2025-08-15 15:16
#
# The "solar" worksheet carries a small lookup column ("grid_cell", column AG,
# rows 4-28) that mirrors which CHE grid-cell each "connecting solar and wind
# to buses..." record (columns R:Y / AC:AG) belongs to. The upstream model
# regenerated that block in a new internal order, so the CHE_xx label that
# belongs with each row changed even though the row positions stayed put.
# Re-point each AG cell at the grid cell that now matches its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$newGridCell = @{
    4  = "CHE_14"
    5  = "CHE_18"
    6  = "CHE_20"
    7  = "CHE_1"
    8  = "CHE_6"
    9  = "CHE_13"
    10 = "CHE_11"
    11 = "CHE_15"
    12 = "CHE_25"
    13 = "CHE_24"
    14 = "CHE_8"
    15 = "CHE_5"
    16 = "CHE_7"
    17 = "CHE_10"
    18 = "CHE_22"
    19 = "CHE_17"
    20 = "CHE_19"
    21 = "CHE_23"
    22 = "CHE_9"
    23 = "CHE_21"
    24 = "CHE_4"
    25 = "CHE_0"
    26 = "CHE_3"
    27 = "CHE_2"
    28 = "CHE_12"
}

foreach ($row in $newGridCell.Keys) {
    $ws.Range("AG$row").Value = $newGridCell[$row]
}
